$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 (Test ID SignUp_Positive_04) with the new "chronic" onboarding data
$ws.Range("B5").Value = "To verify the employee's selection for the 'Lunch and Learnt' chronic program."""
$ws.Range("C5").Value = "Yes"
$ws.Range("I5").Value = "prabhaAutoErsX9791@mailinator.com"
$ws.Range("K5").Value = "PrabhaAutooOBA"
$ws.Range("L5").Value = "automationgxpl"
$ws.Range("U5").Value = "Chronic Disease"

# Update the active cell selection to D6
$ws.Range("D6").Select()
